$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.2265316855114383
$ws.Range("E2").Value = 0.4815058410167694

$ws.Range("C3").Value = 0.6166891672684198
$ws.Range("D3").Value = 0.5721471350418985
$ws.Range("E3").Value = 0.5200661420822144

$ws.Range("D4").Value = 0.7333003895122376
$ws.Range("E4").Value = 0.4010415077209473

$ws.Range("C5").Value = 0.7445664032668693
$ws.Range("D5").Value = 0.5295127493052906
$ws.Range("E5").Value = 0.7511604428291321
